# "Reorder supervision by year" -- swap the two supervision entries in
# rows 6/8 (and their joint-supervision footnote rows 7/9) so the table
# reads in chronological order, then tidy up the trailing blank rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the data rows 6 and 8 (every column, A-E) -----------------------
$cols = @("A", "B", "C", "D", "E")
foreach ($col in $cols) {
    $top = $ws.Range($col + "6").Value2
    $bottom = $ws.Range($col + "8").Value2
    $ws.Range($col + "6").Value2 = $bottom
    $ws.Range($col + "8").Value2 = $top
}

# Row heights travel with their data (row 6 becomes the taller entry,
# row 8 becomes the shorter one).
$h6 = $ws.Rows(6).RowHeight
$h8 = $ws.Rows(8).RowHeight
$ws.Rows(6).RowHeight = $h8
$ws.Rows(8).RowHeight = $h6

# --- Swap the joint-supervision note in rows 7 and 9 (column E only) -----
$e7 = $ws.Range("E7").Value2
$e9 = $ws.Range("E9").Value2
$ws.Range("E7").Value2 = $e9
$ws.Range("E9").Value2 = $e7

# --- Row 10 loses its (empty) styled C cell -------------------------------
$ws.Range("C10").Clear()

# --- Rows 15 and 16 gain an (empty) styled C cell, matching C17/C18 ------
$ws.Range("C17").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Drop the two now-unused trailing blank rows --------------------------
$ws.Rows("19:20").Delete()

# --- Restore the on-screen selection ---------------------------------------
$ws.Range("D15").Select()
